$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-05 Tuesday" "2023-09-06 Wednesday"

Replace-Text "15×35=525" "18×21=378"
Replace-Text "22×93=2046" "34×82=2788"
Replace-Text "93×64=5952" "49×70=3430"
Replace-Text "60×51=3060" "62×16=992"
Replace-Text "52×73=3796" "45×36=1620"

Replace-Text "15×41=615" "87×98=8526"
Replace-Text "61×50=3050" "63×84=5292"
Replace-Text "49×95=4655" "95×18=1710"
Replace-Text "72×88=6336" "40×70=2800"
Replace-Text "86×83=7138" "32×67=2144"

Replace-Text "87×26=2262" "29×68=1972"
Replace-Text "89×20=1780" "16×63=1008"
Replace-Text "29×93=2697" "37×75=2775"
Replace-Text "11×19=209" "49×19=931"
Replace-Text "42×25=1050" "54×60=3240"

Replace-Text "95×14=1330" "36×25=900"
Replace-Text "67×37=2479" "48×44=2112"
Replace-Text "37×59=2183" "49×70=3430"
Replace-Text "97×75=7275" "55×67=3685"
Replace-Text "71×97=6887" "26×87=2262"

Replace-Text "92×63=5796" "52×75=3900"
Replace-Text "42×22=924" "83×67=5561"
Replace-Text "74×91=6734" "92×76=6992"
Replace-Text "28×67=1876" "53×71=3763"
Replace-Text "47×31=1457" "27×54=1458"
